$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" column header in H1, matching the style used by the other
# header cells (copy format from the neighboring "sum" header in G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for the data rows (2-5)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
